$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4618053333333333
$ws.Range("H2").Value = 1.385416
$ws.Range("I2").Value = 0.1241019586052997
$ws.Range("J2").Value = 0.1340391058009509
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4618053333333333
$ws.Range("N2").Value = 1.385416
$ws.Range("O2").Value = 0.1241019586052997
$ws.Range("P2").Value = 0.1340391058009509
$ws.Range("Q2").Value = 0.2132641658951111
$ws.Range("R2").Value = 1.919377493056
$ws.Range("S2").Value = 0.01540129612967151
$ws.Range("T2").Value = 0.0179664818839185
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4618053333333333
$ws.Range("H3").Value = 1.385416
$ws.Range("I3").Value = 0.1241019586052997
$ws.Range("J3").Value = 0.1340391058009509
$ws.Range("O3").Value = 0.6534895031638961
$ws.Range("P3").Value = 0.7058160051525222
$ws.Range("Q3").Value = 1.122995119333333
$ws.Range("R3").Value = 10.106956074
$ws.Range("S3").Value = 0.08109932727064367
$ws.Range("T3").Value = 0.09460694619064339
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4618053333333333
$ws.Range("H4").Value = 1.385416
$ws.Range("I4").Value = 0.1241019586052997
$ws.Range("J4").Value = 0.1340391058009509
$ws.Range("M4").Value = 0.8276215
$ws.Range("N4").Value = 1.655243
$ws.Range("O4").Value = 0.2224085382308043
$ws.Range("P4").Value = 0.160144889046527
$ws.Range("Q4").Value = 0.3822000226813334
$ws.Range("R4").Value = 2.293200136088
$ws.Range("S4").Value = 0.02760133520498449
$ws.Range("T4").Value = 0.02146567772638897
$ws.Range("I5").Value = 0.6534895031638961
$ws.Range("J5").Value = 0.7058160051525222
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4618053333333333
$ws.Range("N5").Value = 1.385416
$ws.Range("O5").Value = 0.1241019586052997
$ws.Range("P5").Value = 0.1340391058009509
$ws.Range("Q5").Value = 1.122995119333333
$ws.Range("R5").Value = 10.106956074
$ws.Range("S5").Value = 0.08109932727064367
$ws.Range("T5").Value = 0.09460694619064339
$ws.Range("I6").Value = 0.6534895031638961
$ws.Range("J6").Value = 0.7058160051525222
$ws.Range("O6").Value = 0.6534895031638961
$ws.Range("P6").Value = 0.7058160051525222
$ws.Range("S6").Value = 0.4270485307453958
$ws.Range("T6").Value = 0.4981762331294652
$ws.Range("I7").Value = 0.6534895031638961
$ws.Range("J7").Value = 0.7058160051525222
$ws.Range("M7").Value = 0.8276215
$ws.Range("N7").Value = 1.655243
$ws.Range("O7").Value = 0.2224085382308043
$ws.Range("P7").Value = 0.160144889046527
$ws.Range("Q7").Value = 2.012568582625
$ws.Range("R7").Value = 12.07541149575
$ws.Range("S7").Value = 0.1453416451478567
$ws.Range("T7").Value = 0.1130328258324136
$ws.Range("G8").Value = 0.8276215
$ws.Range("H8").Value = 1.655243
$ws.Range("I8").Value = 0.2224085382308043
$ws.Range("J8").Value = 0.160144889046527
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4618053333333333
$ws.Range("N8").Value = 1.385416
$ws.Range("O8").Value = 0.1241019586052997
$ws.Range("P8").Value = 0.1340391058009509
$ws.Range("Q8").Value = 0.3822000226813334
$ws.Range("R8").Value = 2.293200136088
$ws.Range("S8").Value = 0.02760133520498449
$ws.Range("T8").Value = 0.02146567772638897
$ws.Range("G9").Value = 0.8276215
$ws.Range("H9").Value = 1.655243
$ws.Range("I9").Value = 0.2224085382308043
$ws.Range("J9").Value = 0.160144889046527
$ws.Range("O9").Value = 0.6534895031638961
$ws.Range("P9").Value = 0.7058160051525222
$ws.Range("Q9").Value = 2.012568582625
$ws.Range("R9").Value = 12.07541149575
$ws.Range("S9").Value = 0.1453416451478567
$ws.Range("T9").Value = 0.1130328258324136
$ws.Range("G10").Value = 0.8276215
$ws.Range("H10").Value = 1.655243
$ws.Range("I10").Value = 0.2224085382308043
$ws.Range("J10").Value = 0.160144889046527
$ws.Range("M10").Value = 0.8276215
$ws.Range("N10").Value = 1.655243
$ws.Range("O10").Value = 0.2224085382308043
$ws.Range("P10").Value = 0.160144889046527
$ws.Range("Q10").Value = 0.68495734726225
$ws.Range("R10").Value = 2.739829389049
$ws.Range("S10").Value = 0.02760133520498449
$ws.Range("T10").Value = 0.02146567772638897
